$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (hunk 0)
$ws.Range("H12").Value = 470.25
$ws.Range("I12").Value = 508.2
$ws.Range("J12").Value = 407
$ws.Range("K12").Value = 508.2
$ws.Range("L12").Value = 407
$ws.Range("M12").Value = -338.2
$ws.Range("N12").Value = -747

# Row 16 (hunk 1)
$ws.Range("H16").Value = 27499.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 27499.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 27499.5
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -27959.5

# Row 17 (hunk 2)
$ws.Range("H17").Value = 6191.125
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 6191.125
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 18573.375
$ws.Range("N17").Value = -18909.375

# Row 107 (hunk 3)
$ws.Range("H107").Value = 1700.5454
$ws.Range("I107").Value = 1462.5
$ws.Range("J107").Value = 2335.3333
$ws.Range("K107").Value = 1462.5
$ws.Range("L107").Value = 2335.3333
$ws.Range("M107").Value = 457.5
$ws.Range("N107").Value = -6175.3333

# Row 113 (hunk 4)
$ws.Range("H113").Value = 1854647.1
$ws.Range("I113").Value = 5558033
$ws.Range("J113").Value = 2954.1667
$ws.Range("K113").Value = 5558033
$ws.Range("L113").Value = 2954.1667
$ws.Range("M113").Value = -5554779
$ws.Range("N113").Value = -9462.1667

# Row 131 (hunk 5)
$ws.Range("H131").Value = 2082.5
$ws.Range("I131").Value = 1899
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 5697
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = -657
$ws.Range("N131").Value = -19080

# Row 137 (hunk 6)
$ws.Range("H137").Value = 2089.25
$ws.Range("I137").Value = 1826.4546
$ws.Range("J137").Value = 2410.4443
$ws.Range("K137").Value = 5479.3638
$ws.Range("L137").Value = 7231.3329
$ws.Range("M137").Value = -2929.3638
$ws.Range("N137").Value = -12331.3329

# Row 138 (hunk 7)
$ws.Range("H138").Value = 2393.7297
$ws.Range("I138").Value = 2228.2222
$ws.Range("J138").Value = 2550.5264
$ws.Range("K138").Value = 6684.6666
$ws.Range("L138").Value = 7651.5792
$ws.Range("M138").Value = -1544.6666
$ws.Range("N138").Value = -17931.5792

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 8)
$ws.Range("H32").Value = 11980.641
$ws.Range("I32").Value = 7543.7827
$ws.Range("J32").Value = 18358.625
$ws.Range("K32").Value = 7543.7827
$ws.Range("L32").Value = 18358.625
$ws.Range("M32").Value = -7256.7827
$ws.Range("N32").Value = -18932.625

# Row 45 (hunk 9)
$ws.Range("H45").Value = 3321.5454
$ws.Range("I45").Value = 2692.375
$ws.Range("J45").Value = 4999.3335
$ws.Range("K45").Value = 2692.375
$ws.Range("L45").Value = 4999.3335
$ws.Range("M45").Value = -2315.375
$ws.Range("N45").Value = -5753.3335

# Row 120 (hunk 10)
$ws.Range("H120").Value = 90499.5
$ws.Range("I120").Value = 81000
$ws.Range("J120").Value = 99999
$ws.Range("K120").Value = 81000
$ws.Range("L120").Value = 99999
$ws.Range("M120").Value = -76162
$ws.Range("N120").Value = -109675

# Row 139 (hunk 11)
$ws.Range("H139").Value = 97762.14
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 101436.16
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 101436.16
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -111716.16

$ws = $wb.Worksheets.Item("BSM")
# Row 80 (hunk 12)
$ws.Range("H80").Value = 303.36365
$ws.Range("I80").Value = 50
$ws.Range("J80").Value = 328.7
$ws.Range("K80").Value = 50
$ws.Range("L80").Value = 328.7
$ws.Range("M80").Value = 948
$ws.Range("N80").Value = -2324.7

# Row 83 (hunk 13)
$ws.Range("H83").Value = 303.36365
$ws.Range("I83").Value = 50
$ws.Range("J83").Value = 328.7
$ws.Range("K83").Value = 250
$ws.Range("L83").Value = 1643.5
$ws.Range("M83").Value = 4742
$ws.Range("N83").Value = -11627.5

# Row 99 (hunk 14)
$ws.Range("H99").Value = 896.3333
$ws.Range("I99").Value = 500
$ws.Range("J99").Value = 1292.6666
$ws.Range("K99").Value = 500
$ws.Range("L99").Value = 1292.6666
$ws.Range("M99").Value = 998
$ws.Range("N99").Value = -4288.6666

# Row 107 (hunk 15)
$ws.Range("H107").Value = 1642.4
$ws.Range("I107").Value = 1194
$ws.Range("J107").Value = 2315
$ws.Range("K107").Value = 1194
$ws.Range("L107").Value = 2315
$ws.Range("M107").Value = 726
$ws.Range("N107").Value = -6155

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 16)
$ws.Range("H31").Value = 8275.465
$ws.Range("I31").Value = 4861.185
$ws.Range("J31").Value = 14037.0625
$ws.Range("K31").Value = 4861.185
$ws.Range("L31").Value = 14037.0625
$ws.Range("M31").Value = -4566.185
$ws.Range("N31").Value = -14627.0625

# Row 34 (hunk 17)
$ws.Range("H34").Value = 8275.465
$ws.Range("I34").Value = 4861.185
$ws.Range("J34").Value = 14037.0625
$ws.Range("K34").Value = 4861.185
$ws.Range("L34").Value = 14037.0625
$ws.Range("M34").Value = -4659.185
$ws.Range("N34").Value = -14441.0625

# Row 58 (hunk 18)
$ws.Range("H58").Value = 3006.3809
$ws.Range("I58").Value = 2794.4
$ws.Range("J58").Value = 3199.0908
$ws.Range("K58").Value = 2794.4
$ws.Range("L58").Value = 3199.0908
$ws.Range("M58").Value = -2591.4
$ws.Range("N58").Value = -3605.0908

# Row 92 (hunk 19)
$ws.Range("H92").Value = 41400
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 41400
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 41400
$ws.Range("N92").Value = -46392

# Row 107 (hunk 20)
$ws.Range("H107").Value = 4291.154
$ws.Range("I107").Value = 4253.8184
$ws.Range("J107").Value = 4496.5
$ws.Range("K107").Value = 4253.8184
$ws.Range("L107").Value = 4496.5
$ws.Range("M107").Value = -2333.8184
$ws.Range("N107").Value = -8336.5

# Row 136 (hunk 21)
$ws.Range("H136").Value = 3006.3809
$ws.Range("I136").Value = 2794.4
$ws.Range("J136").Value = 3199.0908
$ws.Range("K136").Value = 8383.200000000001
$ws.Range("L136").Value = 9597.2724
$ws.Range("M136").Value = -5833.200000000001
$ws.Range("N136").Value = -14697.2724

$ws = $wb.Worksheets.Item("CUL")
# Row 104 (hunk 22)
$ws.Range("H104").Value = 2220
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2220
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 6660
$ws.Range("N104").Value = -11902

# Row 140 (hunk 23)
$ws.Range("H140").Value = 2000.5927
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 2000.5927
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 6001.7781
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -16361.7781

$ws = $wb.Worksheets.Item("GSM")
# Row 92 (hunk 24)
$ws.Range("H92").Value = 40878.75
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 40878.75
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 40878.75
$ws.Range("N92").Value = -44622.75

# Row 102 (hunk 25)
$ws.Range("H102").Value = 1829.8948
$ws.Range("I102").Value = 1184.7931
$ws.Range("J102").Value = 3908.5557
$ws.Range("K102").Value = 1184.7931
$ws.Range("L102").Value = 3908.5557
$ws.Range("M102").Value = 437.2068999999999
$ws.Range("N102").Value = -7152.5557

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (hunk 26)
$ws.Range("H46").Value = 3067.2778
$ws.Range("I46").Value = 954.8
$ws.Range("J46").Value = 3879.7693
$ws.Range("K46").Value = 954.8
$ws.Range("L46").Value = 3879.7693
$ws.Range("M46").Value = -766.8
$ws.Range("N46").Value = -4255.7693

# Row 53 (hunk 27)
$ws.Range("H53").Value = 7500
$ws.Range("I53").Value = 7500
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 7500
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -6982
$ws.Range("N53").ClearContents()

# Row 55 (hunk 28)
$ws.Range("H55").Value = 1996.4445
$ws.Range("I55").Value = 247.5
$ws.Range("J55").Value = 2496.1428
$ws.Range("K55").Value = 247.5
$ws.Range("L55").Value = 2496.1428
$ws.Range("M55").Value = -74.5
$ws.Range("N55").Value = -2842.1428

# Row 61 (hunk 29)
$ws.Range("H61").Value = 4858
$ws.Range("I61").Value = 2033.1111
$ws.Range("J61").Value = 13332.667
$ws.Range("K61").Value = 2033.1111
$ws.Range("L61").Value = 13332.667
$ws.Range("M61").Value = -1831.1111
$ws.Range("N61").Value = -13736.667

# Row 113 (hunk 30)
$ws.Range("H113").Value = 4858
$ws.Range("I113").Value = 2033.1111
$ws.Range("J113").Value = 13332.667
$ws.Range("K113").Value = 2033.1111
$ws.Range("L113").Value = 13332.667
$ws.Range("M113").Value = 136.8888999999999
$ws.Range("N113").Value = -17672.667

# Row 122 (hunk 31)
$ws.Range("H122").Value = 5494.8667
$ws.Range("I122").Value = 4774.909
$ws.Range("J122").Value = 7474.75
$ws.Range("K122").Value = 14324.727
$ws.Range("L122").Value = 22424.25
$ws.Range("M122").Value = -11874.727
$ws.Range("N122").Value = -27324.25

# Row 132 (hunk 32)
$ws.Range("H132").Value = 9319.467000000001
$ws.Range("I132").Value = 8216.566000000001
$ws.Range("J132").Value = 11525.267
$ws.Range("K132").Value = 24649.698
$ws.Range("L132").Value = 34575.801
$ws.Range("M132").Value = -22119.698
$ws.Range("N132").Value = -39635.801

# Row 136 (hunk 33)
$ws.Range("H136").Value = 6689.647
$ws.Range("I136").Value = 5872.3
$ws.Range("J136").Value = 8960.056
$ws.Range("K136").Value = 17616.9
$ws.Range("L136").Value = 26880.168
$ws.Range("M136").Value = -15066.9
$ws.Range("N136").Value = -31980.168

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (hunk 34)
$ws.Range("H107").Value = 1149.08
$ws.Range("I107").Value = 1193.1212
$ws.Range("J107").Value = 1063.5883
$ws.Range("K107").Value = 3579.3636
$ws.Range("L107").Value = 3190.7649
$ws.Range("M107").Value = -1659.3636
$ws.Range("N107").Value = -7030.7649

# Row 113 (hunk 35)
$ws.Range("H113").Value = 700.75
$ws.Range("I113").Value = 761.8333
$ws.Range("J113").Value = 517.5
$ws.Range("K113").Value = 2285.4999
$ws.Range("L113").Value = 1552.5
$ws.Range("M113").Value = -115.4998999999998
$ws.Range("N113").Value = -5892.5
